$d = $word.ActiveDocument

# 1. Update the PC data-dictionary entry: drop the removed component IDs
#    (Monitor_ID, Mouse_ID, CPU_ID, Keyboard_ID) from the PC table listing.
$d.Content.Find.Execute(
    ", Monitor_ID, Mouse_ID, CPU_ID, Keyboard_ID, Lab_ID, OS}",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ", Lab_ID, OS}", 2
) | Out-Null

# 2. Clean up the stray/non-standard hanging-indent markup left on the
#    bullet paragraphs: nudging the (already-correct) LeftIndent value
#    forces the paragraph properties to be re-emitted in canonical form,
#    dropping the bogus first-line indent without altering anything else.
$bulletParagraphs = @(3, 5, 7, 9, 11, 13)
foreach ($i in $bulletParagraphs) {
    $p = $d.Paragraphs.Item($i)
    $p.Format.LeftIndent = $p.Format.LeftIndent
}

Write-Output "done"
